$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A14").Value = 17000011
$ws.Range("H14").Value = "t11"
$ws.Range("B14").Value = "二十一"
$ws.Range("C14").Value = 60
$ws.Range("D14").Value = 75
$ws.Range("E14").Value = 90
$ws.Range("F14").Value = 1110
$ws.Range("G14").Value = "GameButton10"

$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A3:H14"))

$ws.Range("D14").Select()
